# Fruta / hortaliza, semanal
# Insert a new daily price record for "Papa" (Rodeo, 1a (guarda lavada)) at row 202,
# pushing the existing rows 202:296 down to 203:297.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 202 downward (Excel default shift = down for an inserted row)
$ws.Rows("202:202").Insert()

# Populate the newly inserted row 202 with the new record
$ws.Range("A202").Value = 5
$ws.Range("B202").Value = "Macroferia Regional de Talca"
$ws.Range("C202").Value = "Maule"
$ws.Range("D202").Value = 44489
$ws.Range("E202").Value = 7
$ws.Range("F202").Value = 100114001
$ws.Range("G202").Value = "Papa"
$ws.Range("H202").Value = "Rodeo"
$ws.Range("I202").Value = "1a (guarda lavada)"
$ws.Range("J202").Value = 1500
$ws.Range("K202").Value = 10000
$ws.Range("L202").Value = 10000
$ws.Range("M202").Value = 10000
$ws.Range("N202").Value = "`$/malla 25 kilos"
$ws.Range("O202").Value = "Región de La Araucanía"
$ws.Range("P202").Value = 400
$ws.Range("Q202").Value = 25
$ws.Range("R202").Value = "Hortaliza"
